$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 3
$ws.Range("I3").Value = 2.45
$ws.Range("AG3").Value = 9.5
$ws.Range("AI3").Value = 19
# Row 5
$ws.Range("G5").Value = 3.5
$ws.Range("I5").Value = 2.5
$ws.Range("R5").Value = 2.25
$ws.Range("S5").Value = 1.57
$ws.Range("U5").Value = 15
$ws.Range("W5").Value = 41
$ws.Range("Z5").Value = 5.5
$ws.Range("AF5").Value = 10
$ws.Range("AI5").Value = 26
# Row 7
$ws.Range("J7").Value = 1.04
$ws.Range("K7").Value = 13
$ws.Range("P7").Value = 1.36
$ws.Range("Q7").Value = 3
$ws.Range("U7").Value = 8.5
$ws.Range("Z7").Value = 13
$ws.Range("AA7").Value = 8
$ws.Range("AB7").Value = 15
$ws.Range("AC7").Value = 41
# Row 12
$ws.Range("G12").Value = 1.75
$ws.Range("H12").Value = 3.2
$ws.Range("I12").Value = 4.9
$ws.Range("L12").Value = 1.47
$ws.Range("M12").Value = 2.32
$ws.Range("N12").Value = 2.37
$ws.Range("O12").Value = 1.45
$ws.Range("P12").Value = 1.55
$ws.Range("Q12").Value = 2.15
$ws.Range("R12").Value = 2.2
$ws.Range("S12").Value = 1.52
$ws.Range("T12").Value = 4.9
$ws.Range("U12").Value = 6.7
$ws.Range("V12").Value = 9
$ws.Range("W12").Value = 13
$ws.Range("X12").Value = 18
$ws.Range("Y12").Value = 45
$ws.Range("Z12").Value = 6.5
$ws.Range("AA12").Value = 6.5
$ws.Range("AB12").Value = 22
$ws.Range("AC12").Value = 150
$ws.Range("AE12").Value = 10
$ws.Range("AF12").Value = 27
$ws.Range("AG12").Value = 17.5
$ws.Range("AH12").Value = 100
$ws.Range("AI12").Value = 65
# Row 15
$ws.Range("G15").Value = 4.65
$ws.Range("H15").Value = 3.5
$ws.Range("I15").Value = 1.65
$ws.Range("N15").Value = 1.91
$ws.Range("O15").Value = 1.7
$ws.Range("U15").Value = 21
$ws.Range("W15").Value = 60
$ws.Range("X15").Value = 37
$ws.Range("Y15").Value = 40
$ws.Range("Z15").Value = 9
$ws.Range("AA15").Value = 6.1
$ws.Range("AB15").Value = 14
$ws.Range("AC15").Value = 65
$ws.Range("AE15").Value = 5.4
$ws.Range("AF15").Value = 6.2
$ws.Range("AH15").Value = 10
$ws.Range("AI15").Value = 11.25
$ws.Range("AJ15").Value = 23
# Row 16
$ws.Range("G16").Value = 2.3
$ws.Range("H16").Value = 3.4
$ws.Range("I16").Value = 2.62
$ws.Range("O16").Value = 1.85
$ws.Range("T16").Value = 7.5
$ws.Range("U16").Value = 9.75
$ws.Range("V16").Value = 7.8
$ws.Range("W16").Value = 18.5
$ws.Range("X16").Value = 14.5
$ws.Range("Y16").Value = 21
$ws.Range("Z16").Value = 11
$ws.Range("AA16").Value = 5.9
$ws.Range("AE16").Value = 8
$ws.Range("AF16").Value = 11.5
$ws.Range("AG16").Value = 8.5
$ws.Range("AH16").Value = 23
$ws.Range("AI16").Value = 17
# Row 17
$ws.Range("G17").Value = 2.25
$ws.Range("H17").Value = 3.1
$ws.Range("J17").Value = 1.08
$ws.Range("K17").Value = 8
$ws.Range("L17").Value = 1.4
$ws.Range("M17").Value = 2.75
$ws.Range("N17").Value = 2.25
$ws.Range("O17").Value = 1.62
$ws.Range("R17").Value = 1.95
$ws.Range("S17").Value = 1.8
$ws.Range("T17").Value = 7
$ws.Range("X17").Value = 21
$ws.Range("Y17").Value = 34
$ws.Range("Z17").Value = 7.5
$ws.Range("AD17").Value = 351
$ws.Range("AE17").Value = 9
$ws.Range("AG17").Value = 13
$ws.Range("AH17").Value = 41
# Row 18
$ws.Range("L18").Value = 1.4
$ws.Range("M18").Value = 2.75
# Row 23
$ws.Range("G23").Value = 3.2
$ws.Range("H23").Value = 3.3
$ws.Range("I23").Value = 2.1
$ws.Range("J23").Value = 1.06
$ws.Range("K23").Value = 10
$ws.Range("L23").Value = 1.3
$ws.Range("M23").Value = 3.4
$ws.Range("N23").Value = 2
$ws.Range("P23").Value = 1.4
$ws.Range("Q23").Value = 2.75
$ws.Range("R23").Value = 1.8
$ws.Range("S23").Value = 1.91
$ws.Range("U23").Value = 17
$ws.Range("V23").Value = 12
$ws.Range("W23").Value = 34
$ws.Range("X23").Value = 26
$ws.Range("Y23").Value = 34
$ws.Range("Z23").Value = 10
$ws.Range("AA23").Value = 6.5
$ws.Range("AB23").Value = 15
$ws.Range("AC23").Value = 51
$ws.Range("AD23").Value = 251
$ws.Range("AE23").Value = 7.5
$ws.Range("AF23").Value = 10
$ws.Range("AG23").Value = 9
$ws.Range("AH23").Value = 19
$ws.Range("AI23").Value = 17
$ws.Range("AJ23").Value = 29
# Row 24
$ws.Range("G24").Value = 3.3
$ws.Range("I24").Value = 2
$ws.Range("J24").Value = 1.05
$ws.Range("K24").Value = 11
$ws.Range("L24").Value = 1.25
$ws.Range("M24").Value = 3.75
$ws.Range("O24").Value = 1.95
$ws.Range("P24").Value = 1.36
$ws.Range("Q24").Value = 3
$ws.Range("R24").Value = 1.73
$ws.Range("S24").Value = 2
$ws.Range("T24").Value = 11
$ws.Range("U24").Value = 19
$ws.Range("V24").Value = 12
$ws.Range("W24").Value = 41
$ws.Range("X24").Value = 26
$ws.Range("Y24").Value = 34
$ws.Range("Z24").Value = 11
$ws.Range("AA24").Value = 6.5
$ws.Range("AB24").Value = 13
$ws.Range("AC24").Value = 41
$ws.Range("AD24").Value = 201
$ws.Range("AE24").Value = 8
$ws.Range("AF24").Value = 10
$ws.Range("AH24").Value = 19
$ws.Range("AJ24").Value = 26
# Row 26
$ws.Range("H26").Value = 4
$ws.Range("J26").Value = 1.03
$ws.Range("K26").Value = 10.5
$ws.Range("L26").Value = 1.17
$ws.Range("M26").Value = 4.5
$ws.Range("N26").Value = 1.53
$ws.Range("O26").Value = 2.38
$ws.Range("P26").Value = 1.25
$ws.Range("Q26").Value = 3.75
$ws.Range("R26").Value = 1.57
$ws.Range("S26").Value = 2.25
$ws.Range("T26").Value = 17
$ws.Range("U26").Value = 23
$ws.Range("Z26").Value = 17
$ws.Range("AB26").Value = 13
$ws.Range("AD26").Value = 126
$ws.Range("AE26").Value = 10
# Row 33
$ws.Range("G33").Value = 6
$ws.Range("I33").Value = 1.5
$ws.Range("L33").Value = 1.25
$ws.Range("M33").Value = 3.75
$ws.Range("N33").Value = 1.88
$ws.Range("O33").Value = 1.98
$ws.Range("R33").Value = 1.95
$ws.Range("S33").Value = 1.8
$ws.Range("AC33").Value = 51
$ws.Range("AD33").Value = 351
$ws.Range("AJ33").Value = 26
